$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.929.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2561"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06363"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.282"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.647.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5438"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7812"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.959.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "196.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.432"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.932"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.044"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.896"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1166"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.66%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.235"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04958"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.253"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.178"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.536"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.363"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8920"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.588"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.133.61"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5435"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.11%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.546"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₈128"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.586"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8155"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.776.75"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4541"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.57"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05073"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.50%  "
